# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets, which contain duplicate data tables.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 11642
    "F3"  = 11212
    "F6"  = 1013
    "F7"  = 120
    "F8"  = 69
    "F9"  = 42
    "F10" = 44
    "F11" = 10701
    "F12" = 4138
    "F13" = 14
    "F17" = 1050
    "F18" = 47
    "F20" = 438
    "F21" = 11122
    "F22" = 10888
    "F24" = 25
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
